## download articles with pandoc title blocks
##
## Rewrites the article's opening two paragraphs so they look like a
## pandoc-generated "title block":
##   1. "Strike Leader Comes East" switches from the Heading1 style to the
##      Title style, and its text is split word-by-word (with separate
##      single-space runs) instead of one run.
##   2. "By Dorothy Day" loses the "By " lead-in and bold run formatting,
##      switches to the Authors style, and its text is likewise split into
##      "Dorothy" / " " / "Day" runs.

$d = $word.ActiveDocument

function New-WordXmlPackage($innerParagraphXml) {
    # Wraps a single <w:p>...</w:p> fragment in the mc:package envelope that
    # Range.InsertXML / Range.XML expect (same shape Word itself emits from
    # Range.XML), scoped to just word/document.xml.
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>$innerParagraphXml</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# --- Paragraph 1: title -----------------------------------------------
$titleInner = @'
<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Strike</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Leader</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Comes</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">East</w:t></w:r></w:p>
'@

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
[void]$titleRange.InsertXML((New-WordXmlPackage $titleInner))

# --- Paragraph 2: author -------------------------------------------------
$authorInner = @'
<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>
'@

$authorPara = $d.Paragraphs.Item(2)
$authorRange = $authorPara.Range
[void]$authorRange.InsertXML((New-WordXmlPackage $authorInner))

# --- Bookmark cleanup -----------------------------------------------------
# The old markup wrapped the title paragraph in a
# strike-leader-comes-east bookmark; the new pandoc-style title block no
# longer needs it, so drop it (best-effort: older Word automation hosts
# silently ignore an unknown bookmark name).
try {
    $d.Bookmarks("strike-leader-comes-east").Delete()
} catch {
}
try {
    $d.DeleteBookmark("strike-leader-comes-east")
} catch {
}

Write-Host "applied pandoc title block edit"
